# Applies the per-cell coin price / volume refresh described by the commit diff.
# Columns D (Price) and E (Volume 1h) hold plain text in the source workbook
# (not numbers/percentages), so for column D -- where many values look numeric --
# we force the cell to Text format before writing, then clear formats again so the
# cell keeps no explicit style (matching the original file's unstyled data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "59.120.34"
$ws.Range("E2").Value = "  -0.03%  "
Set-TextValue "D3" "2.515.38"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.16%  "
Set-TextValue "D5" "533.81"
$ws.Range("E5").Value = "  -0.45%  "
Set-TextValue "D6" "139.92"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("E7").Value = "  +0.51%  "
Set-TextValue "D8" "0.563"
$ws.Range("E8").Value = "  -1.88%  "
Set-TextValue "D9" "2.520.11"
$ws.Range("E9").Value = "  -0.53%  "
Set-TextValue "D10" "0.0989"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("E11").Value = "  +1.57%  "
Set-TextValue "D12" "5.41"
$ws.Range("E12").Value = "  -2.61%  "
Set-TextValue "D13" "0.352"
$ws.Range("E13").Value = "  -0.35%  "
Set-TextValue "D14" "2.964.89"
$ws.Range("E14").Value = "  +1.05%  "
Set-TextValue "D15" "23.19"
$ws.Range("E15").Value = "  -3.44%  "
Set-TextValue "D16" "59.074.49"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +0.39%  "
Set-TextValue "D18" "2.515.75"
$ws.Range("E18").Value = "  +0.35%  "
Set-TextValue "D19" "10.93"
$ws.Range("E19").Value = "  -3.81%  "
Set-TextValue "D20" "4.22"
$ws.Range("E20").Value = "  -1.70%  "
Set-TextValue "D21" "320.00"
$ws.Range("E21").Value = "  -1.23%  "
Set-TextValue "D22" "1.00"
$ws.Range("E22").Value = "  +0.10%  "
Set-TextValue "D23" "5.77"
$ws.Range("E23").Value = "  +0.12%  "
Set-TextValue "D24" "62.24"
$ws.Range("E24").Value = "  +2.16%  "
Set-TextValue "D25" "0.419"
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("E26").Value = "  +1.16%  "
Set-TextValue "D27" "1.00"
$ws.Range("E27").Value = "  +0.74%  "
Set-TextValue "D28" "7.78"
$ws.Range("E28").Value = "  -0.03%  "
Set-TextValue "D29" "6.77"
$ws.Range("E29").Value = "  -2.38%  "
Set-TextValue "D30" "0.0₃0765"
$ws.Range("E30").Value = "  -2.56%  "
Set-TextValue "D31" "1.78"
$ws.Range("E31").Value = "  -0.36%  "
Set-TextValue "D32" "163.97"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("E34").Value = "  -9.25%  "
$ws.Range("E35").Value = "  -0.72%  "
Set-TextValue "D36" "18.38"
$ws.Range("E36").Value = "  -0.96%  "
Set-TextValue "D37" "4.25"
$ws.Range("E37").Value = "  -4.25%  "
Set-TextValue "D38" "1.57"
$ws.Range("E38").Value = "  -2.68%  "
Set-TextValue "D39" "36.81"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D40" "5.38"
$ws.Range("E40").Value = "  -8.82%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D41" "3.62"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D42" "287.86"
$ws.Range("E42").Value = "  -7.37%  "
Set-TextValue "D43" "0.804"
$ws.Range("E43").Value = "  -1.92%  "
Set-TextValue "D44" "0.998"
$ws.Range("E44").Value = "  +0.39%  "
Set-TextValue "D45" "10.86"
$ws.Range("E45").Value = "  +0.86%  "
Set-TextValue "D46" "0.595"
$ws.Range("E46").Value = "  -0.58%  "
Set-TextValue "D47" "123.92"
$ws.Range("E47").Value = "  -1.04%  "
Set-TextValue "D48" "0.0927"
$ws.Range("E48").Value = "  -0.35%  "
Set-TextValue "D49" "18.51"
$ws.Range("E49").Value = "  -1.15%  "
Set-TextValue "D50" "0.0508"
$ws.Range("E50").Value = "  -1.95%  "
Set-TextValue "D51" "0.0222"
$ws.Range("E51").Value = "  -2.60%  "
